$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.929.60'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '3.463.43'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.44'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.07'
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.584'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("E9").Value = '  -4.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.06'
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("E11").Value = '  -3.52%  '
$ws.Range("D12").Value = '4.062.35'
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '31.11'
$ws.Range("E13").Value = '  +8.09%  '
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '66.917.95'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("E16").Value = '  -4.44%  '
$ws.Range("D17").Value = '3.460.88'
$ws.Range("E17").Value = '  -2.30%  '
$ws.Range("E18").Value = '  -2.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.21'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '384.06'
$ws.Range("E20").Value = '  -3.55%  '
$ws.Range("E21").Value = '  -2.38%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.46'
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("E26").Value = '  -2.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.21'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  -4.03%  '
$ws.Range("E31").Value = '  -4.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.03'
$ws.Range("E32").Value = '  -2.81%  '
$ws.Range("E33").Value = '  -3.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.19'
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.92'
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.866'
$ws.Range("E37").Value = '  -3.43%  '
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.88'
$ws.Range("E39").Value = '  -1.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.93'
$ws.Range("E41").Value = '  -3.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.00'
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("D43").Value = '2.781.13'
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("E44").Value = '  -4.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.55'
$ws.Range("E45").Value = '  -2.96%  '
$ws.Range("E46").Value = '  -2.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0297'
$ws.Range("E47").Value = '  -4.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '335.33'
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("E49").Value = '  -4.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '32.94'
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("E51").Value = '  -3.46%  '
